$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: financial period headers ---
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/10"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/10"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/10"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/10"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/10"

# --- Row 9: publish dates ---
$ws.Range("D9").Value = "1399-04-02 (10)"
$ws.Range("E9").Value = "1400-02-26 (8)"
$ws.Range("F9").Value = "1401-02-19 (9)"
$ws.Range("G9").Value = "1402-02-13 (9)"
$ws.Range("H9").Value = "1402-02-29 (3)"

# --- Balance sheet data rows ---
# Row 12
$ws.Range("D12").Value = 39000
$ws.Range("E12").Value = 49811
$ws.Range("F12").Value = 290638
$ws.Range("G12").Value = 287150
$ws.Range("H12").Value = 586153

# Row 13
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 37000
$ws.Range("G13").Value = 1308000
$ws.Range("H13").Value = 1984062

# Row 14
$ws.Range("D14").Value = 380901
$ws.Range("E14").Value = 1106795
$ws.Range("F14").Value = 1116836
$ws.Range("G14").Value = 1227369
$ws.Range("H14").Value = 1058983

# Row 15
$ws.Range("D15").Value = 692723
$ws.Range("E15").Value = 597355
$ws.Range("F15").Value = 848935
$ws.Range("G15").Value = 1916387
$ws.Range("H15").Value = 2781692

# Row 16
$ws.Range("D16").Value = 133599
$ws.Range("E16").Value = 119528
$ws.Range("F16").Value = 283959
$ws.Range("G16").Value = 174125
$ws.Range("H16").Value = 394478

# Row 17
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 40995

# Row 18
$ws.Range("D18").Value = 1246223
$ws.Range("E18").Value = 1873489
$ws.Range("F18").Value = 2577368
$ws.Range("G18").Value = 4913031
$ws.Range("H18").Value = 6846363

# Row 19
$ws.Range("D19").Value = 2373
$ws.Range("E19").Value = 1669
$ws.Range("F19").Value = 2438
$ws.Range("G19").Value = 2068
$ws.Range("H19").Value = 101577

# Row 20
$ws.Range("D20").Value = 483210
$ws.Range("E20").Value = 440561
$ws.Range("F20").Value = 441350
$ws.Range("G20").Value = 441350
$ws.Range("H20").Value = 441349

# Row 22
$ws.Range("D22").Value = 1065550
$ws.Range("E22").Value = 1024897
$ws.Range("F22").Value = 1484618
$ws.Range("G22").Value = 1650971
$ws.Range("H22").Value = 1722104

# Row 23
$ws.Range("D23").Value = 33559
$ws.Range("E23").Value = 33413
$ws.Range("F23").Value = 33267
$ws.Range("G23").Value = 41198
$ws.Range("H23").Value = 50342

# Row 25
$ws.Range("D25").Value = 6909
$ws.Range("E25").Value = 6909
$ws.Range("F25").Value = 50027
$ws.Range("G25").Value = 78017
$ws.Range("H25").Value = 37036

# Row 26
$ws.Range("D26").Value = 1591601
$ws.Range("E26").Value = 1507449
$ws.Range("F26").Value = 2011700
$ws.Range("G26").Value = 2213604
$ws.Range("H26").Value = 2352408

# Row 27
$ws.Range("D27").Value = 2837824
$ws.Range("E27").Value = 3380938
$ws.Range("F27").Value = 4589068
$ws.Range("G27").Value = 7126635
$ws.Range("H27").Value = 9198771

# Row 29
$ws.Range("D29").Value = 590165
$ws.Range("E29").Value = 555964
$ws.Range("F29").Value = 783601
$ws.Range("G29").Value = 1365260
$ws.Range("H29").Value = 1441599

# Row 31
$ws.Range("D31").Value = 38208
$ws.Range("E31").Value = 29032
$ws.Range("F31").Value = 261612
$ws.Range("G31").Value = 212517
$ws.Range("H31").Value = 208067

# Row 32
$ws.Range("D32").Value = 29718
$ws.Range("E32").Value = 141541
$ws.Range("F32").Value = 249548
$ws.Range("G32").Value = 500612
$ws.Range("H32").Value = 556654

# Row 33
$ws.Range("D33").Value = 54605
$ws.Range("E33").Value = 8387
$ws.Range("F33").Value = 47124
$ws.Range("G33").Value = 295394
$ws.Range("H33").Value = 321427

# Row 34
$ws.Range("D34").Value = 772482
$ws.Range("E34").Value = 577636
$ws.Range("F34").Value = 76321
$ws.Range("G34").Value = 212155
$ws.Range("H34").Value = 0

# Row 37
$ws.Range("D37").Value = 1485178
$ws.Range("E37").Value = 1312560
$ws.Range("F37").Value = 1418206
$ws.Range("G37").Value = 2585938
$ws.Range("H37").Value = 2527747

# Row 38
$ws.Range("D38").Value = 0
$ws.Range("E38").Value = 0
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 53690
$ws.Range("H38").Value = 148197

# Row 39
$ws.Range("D39").Value = 0
$ws.Range("E39").Value = "-"
$ws.Range("F39").Value = "-"
$ws.Range("G39").Value = "-"
$ws.Range("H39").Value = "-"

# Row 40
$ws.Range("D40").Value = 0
$ws.Range("E40").Value = 79227
$ws.Range("F40").Value = 0
$ws.Range("G40").Value = 0
$ws.Range("H40").Value = 0

# Row 41
$ws.Range("D41").Value = 96740
$ws.Range("E41").Value = 140658
$ws.Range("F41").Value = 147335
$ws.Range("G41").Value = 204877
$ws.Range("H41").Value = 316045

# Row 42
$ws.Range("D42").Value = 96740
$ws.Range("E42").Value = 219885
$ws.Range("F42").Value = 147335
$ws.Range("G42").Value = 258567
$ws.Range("H42").Value = 464242

# Row 43
$ws.Range("D43").Value = 1581918
$ws.Range("E43").Value = 1532445
$ws.Range("F43").Value = 1565541
$ws.Range("G43").Value = 2844505
$ws.Range("H43").Value = 2991989

# Row 45
$ws.Range("D45").Value = 502900
$ws.Range("E45").Value = 502900
$ws.Range("F45").Value = 728789
$ws.Range("G45").Value = 728789
$ws.Range("H45").Value = 728789

# Row 48
$ws.Range("D48").Value = 0
$ws.Range("E48").Value = 0
$ws.Range("F48").Value = -88248
$ws.Range("G48").Value = -84053
$ws.Range("H48").Value = -89016

# Row 49
$ws.Range("D49").Value = "-"
$ws.Range("E49").Value = 0
$ws.Range("F49").Value = 1948
$ws.Range("G49").Value = 0
$ws.Range("H49").Value = 0

# Row 50
$ws.Range("D50").Value = 50290
$ws.Range("E50").Value = 50290
$ws.Range("F50").Value = 72879
$ws.Range("G50").Value = 72879
$ws.Range("H50").Value = 72879

# Row 52
$ws.Range("D52").Value = 0
$ws.Range("E52").Value = "-"
$ws.Range("F52").Value = "-"
$ws.Range("G52").Value = "-"
$ws.Range("H52").Value = "-"

# Row 53
$ws.Range("D53").Value = 208739
$ws.Range("E53").Value = 208739
$ws.Range("F53").Value = 208739
$ws.Range("G53").Value = 626371
$ws.Range("H53").Value = 626371

# Row 54
$ws.Range("D54").Value = 0
$ws.Range("E54").Value = "-"
$ws.Range("F54").Value = "-"
$ws.Range("G54").Value = "-"
$ws.Range("H54").Value = "-"

# Row 56
$ws.Range("D56").Value = 493977
$ws.Range("E56").Value = 1086564
$ws.Range("F56").Value = 2099420
$ws.Range("G56").Value = 2938144
$ws.Range("H56").Value = 4867759

# Row 57
$ws.Range("D57").Value = 1255906
$ws.Range("E57").Value = 1848493
$ws.Range("F57").Value = 3023527
$ws.Range("G57").Value = 4282130
$ws.Range("H57").Value = 6206782

# Row 58
$ws.Range("D58").Value = 2837824
$ws.Range("E58").Value = 3380938
$ws.Range("F58").Value = 4589068
$ws.Range("G58").Value = 7126635
$ws.Range("H58").Value = 9198771
